# Minor edits to baseline docs (i.e., text/wording)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Raw cover" ---
$ws1 = $wb.Worksheets.Item("Raw cover")

# Clear the sample/example data row (row 3), but keep the date format on B3
$ws1.Range("A3:E3").ClearContents()

# Reset the view: scroll back so column A is visible, select I6
$ws1.Activate()
$ws1.Range("A1").Select()
$ws1.Range("I6").Select()

# --- Sheet 2: "Metadata" ---
$ws2 = $wb.Worksheets.Item("Metadata")

# Fix wording of the Ab_ description
$ws2.Range("B16").Value = "Counted abundance of an animal species (*Add further columns as needed)"

# Update selection on Metadata sheet
$ws2.Activate()
$ws2.Range("D19").Select()

# Return focus to the first sheet, matching tabSelected="1" on "Raw cover"
$ws1.Activate()
